$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old A1 row (value 0, bold+bordered style) by deleting it entirely;
# this shifts the former A2 (shared-string question text) up to A1 and clears the dimension/style bloat.
$ws.Rows(1).Delete()

# Replace the cell text with the pretty-printed / reformatted questions payload
$questions = @"
questions = [
    {
        "title": "What is the priority of hooks defined in mixins?",
        "ques_type": 2,
        "options": [
            "They're to be launched before component's hooks.",
            "Equal with own component's hooks.",
            "They're to be launched instead of component's hooks.",
            "They'll be overridden with component's hooks."
        ],
        "score": "They're to be launched before component's hooks."
    },
    {
        "title": "Are the results of a Vuex getter cached?",
        "ques_type": 2,
        "options": [
            "Yes, if it does not return a function, otherwise no",
            "Yes, if it returns a function, otherwise no",
            "Yes",
            "No"
        ],
        "score": "Yes, if it does not return a function, otherwise no"
    },
    {
        "title": "What is the concept of Hydration in Vue?",
        "ques_type": 2,
        "options": [
            "A client-side process during which Vue takes over the static HTML sent by the server and turns it into dynamic DOM that can react to client-side data changes.",
            "A client-side process of mounting a server-side rendered page and executing remaining lifecycle hooks.",
            "A server-side process of preparing a state snapshot that corresponds to a rendered page.",
            "A server-side process of injection a client-dependent state into an initial state of an application."
        ],
        "score": "A client-side process during which Vue takes over the static HTML sent by the server and turns it into dynamic DOM that can react to client-side data changes."
    },
    {
        "title": "Where can Vue filters be applied?",
        "ques_type": 2,
        "options": [
            "Mustache interpolations and 'v-bind' expressions",
            "Mustache interpolations",
            "'v-bind' expressions",
            "'v-bind' and 'v-html' expressions"
        ],
        "score": "Mustache interpolations and 'v-bind' expressions"
    }
]
"@
$ws.Range("A1").Value = $questions

# Avoid a stray custom row height from the multi-line content (matches original default-height layout)
$ws.Rows(1).AutoFit()
